$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume snapshot refresh (GitHub Actions bot).
# Column D ("Price") holds free-form text (e.g. "62.581.55", "1.00",
# "0.0₃0997") rather than real numbers, so every write there first forces
# a Text format and restores the original "Normal" style afterwards -
# this keeps Excel from "helpfully" re-parsing the string as a number
# (which would silently drop meaningful trailing/duplicate-looking zeros)
# while leaving the cell style byte-for-byte as it was.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.581.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.097.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.21%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -10.98%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.581'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -9.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.093.39'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.14%  '
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.114'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.13%  '
$ws.Range('E12').Value = '  -6.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.633.51'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.27%  '
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.651.39'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '24.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.097.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.37%  '
$ws.Range('E18').Value = '  -7.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '396.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.01'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.85%  '
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('E24').Value = '  -1.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.194'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.475'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0997'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -12.91%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.54'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.15%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('E32').Value = '  -8.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.15'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '153.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.08'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.697.58'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.75%  '
$ws.Range('E40').Value = '  -8.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.13'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -12.06%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.20%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.97'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.686'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.56%  '
$ws.Range('E45').Value = '  -8.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -12.39%  '
$ws.Range('E47').Value = '  -6.74%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.97%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '278.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -11.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0970'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.55%  '
